# "Add files via upload" -- refresh the cached "today" date shown by the
# Date Placeholder fields (slide master + every slide layout) from
# 9/19/2024 to 9/22/2024, and tighten up the copy on the "BENEFITS FOR
# USERS" slide (slide 4).

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.TextFrame.TextRange.Text -eq "9/19/2024") {
            $sh.TextFrame.TextRange.Text = "9/22/2024"
        }
    }
}

# Slide master's own Date Placeholder.
Update-DateField $p.SlideMaster.Shapes

# Every slide layout has its own copy of the Date Placeholder too.
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# Slide 4 ("BENEFITS FOR USERS"): shrink the benefits textbox to fit its
# now-shorter copy, and drop the trailing clause from the sentence.
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $sh = $s4.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and `
        $sh.TextFrame.TextRange.Text -like "*By utilizing Diet Tracker*") {
        $sh.Height = 181.75779527559055
        $sh.TextFrame.TextRange.Text = "By utilizing Diet Tracker, users can simplify their meal tracking, adopt healthier eating habits, minimize food waste, and cultivate a mindset of sustainability."
        break
    }
}
